# Update scripts with new TPM values.
# The "ECs" sending-cluster row is removed, and the remaining
# "Resolving-Mac" row is recalculated with the new TPM-derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 2 (Sending cluster = "ECs"); this shifts the old
# row 3 ("Resolving-Mac") up into row 2 and drops the now-unused "ECs"
# shared string automatically on save.
$ws.Rows(2).Delete()

# The "Target cluster" column still pointed at the old "ECs" string;
# it should read "Resolving-Mac" now (same as column A / the sender).
$ws.Range("D2").Value = "Resolving-Mac"

# Refresh the recalculated NATMI metrics for the remaining row.
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.501936999999999
$ws.Range("H2").Value = 13.505811
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.006530999999999999
$ws.Range("N2").Value = 0.019593
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.02940215054699999
$ws.Range("R2").Value = 0.2646193549229999
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
